# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 18920
$wsExhibition.Range("F29").Value = 6028
$wsExhibition.Range("F35").Value = 5393
$wsExhibition.Range("F36").Value = 2
$wsExhibition.Range("F37").Value = 7

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 18920
$wsAll.Range("F32").Value = 6028
$wsAll.Range("F38").Value = 5393
$wsAll.Range("F39").Value = 2
$wsAll.Range("F40").Value = 7
